$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds date-like text ("2025-09-10"). Force it to Text format
# before assignment so Excel stores it as a literal string instead of
# silently converting it to a date serial number.
$ws.Range("A15").NumberFormat = "@"

$ws.Range("A15").Value = "2025-09-10"
$ws.Range("B15").Value = "15:19:43"
$ws.Range("C15").Value = "1.00 EUR = 1678.2062 ARS"
